$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 545.4
$ws.Range("J2").Value = 596.75
$ws.Range("L2").Value = 596.75
$ws.Range("N2").Value = -822.75

$ws.Range("H17").Value = 2562.5
$ws.Range("J17").Value = 2000
$ws.Range("L17").Value = 6000
$ws.Range("N17").Value = -6336

$ws.Range("H18").Value = 1630.6154
$ws.Range("I18").Value = 1663.4546
$ws.Range("K18").Value = 1663.4546
$ws.Range("M18").Value = -1379.4546

$ws.Range("H28").Value = 519.5714
$ws.Range("I28").Value = 189.5
$ws.Range("K28").Value = 189.5
$ws.Range("M28").Value = 295.5

$ws.Range("H33").Value = 258.25
$ws.Range("I33").Value = 93
$ws.Range("K33").Value = 93
$ws.Range("M33").Value = 136

$ws.Range("H62").Value = 6036.091
$ws.Range("I62").Value = 4479.6
$ws.Range("K62").Value = 4479.6
$ws.Range("M62").Value = -3855.6

$ws.Range("H64").Value = 3200
$ws.Range("J64").Value = 3200
$ws.Range("L64").Value = 3200
$ws.Range("N64").Value = -3696

$ws.Range("H65").Value = 6036.091
$ws.Range("I65").Value = 4479.6
$ws.Range("K65").Value = 22398
$ws.Range("M65").Value = -19278

$ws.Range("H67").Value = 3200
$ws.Range("J67").Value = 3200
$ws.Range("L67").Value = 3200
$ws.Range("N67").Value = -4916

$ws.Range("H97").Value = 2300
$ws.Range("J97").Value = 2300
$ws.Range("L97").Value = 6900
$ws.Range("N97").Value = -7892

$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()

$ws.Range("H116").Value = 5947.5713
$ws.Range("I116").Value = 5588.222
$ws.Range("J116").Value = 6594.4
$ws.Range("K116").Value = 5588.222
$ws.Range("L116").Value = 6594.4
$ws.Range("M116").Value = -2146.222
$ws.Range("N116").Value = -13478.4

$ws.Range("H135").Value = 775
$ws.Range("I135").Value = 596.9167
$ws.Range("K135").Value = 5372.2503
$ws.Range("M135").Value = -2837.2503

$ws.Range("H137").Value = 1958.125
$ws.Range("I137").Value = 986.73334
$ws.Range("J137").Value = 3577.111
$ws.Range("K137").Value = 2960.20002
$ws.Range("L137").Value = 10731.333
$ws.Range("M137").Value = -410.2000200000002
$ws.Range("N137").Value = -15831.333

$ws.Range("H138").Value = 4016.016
$ws.Range("I138").Value = 1395.8
$ws.Range("K138").Value = 4187.4
$ws.Range("M138").Value = 952.6000000000004

$ws.Range("H141").Value = 1815.1428
$ws.Range("I141").Value = 1647.0769
$ws.Range("K141").Value = 4941.2307
$ws.Range("M141").Value = 238.7692999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 562
$ws.Range("I4").Value = 86
$ws.Range("K4").Value = 86
$ws.Range("M4").Value = 30

$ws.Range("H32").Value = 16661.223
$ws.Range("I32").Value = 3133.5264
$ws.Range("K32").Value = 3133.5264
$ws.Range("M32").Value = -2846.5264

$ws.Range("H61").Value = 3305.6667
$ws.Range("I61").Value = 3305.6667
$ws.Range("K61").Value = 3305.6667
$ws.Range("M61").Value = -3093.6667

$ws.Range("H74").Value = 2561
$ws.Range("I74").Value = 1155
$ws.Range("K74").Value = 1155
$ws.Range("M74").Value = -281

$ws.Range("H77").Value = 2561
$ws.Range("I77").Value = 1155
$ws.Range("K77").Value = 5775
$ws.Range("M77").Value = -1407

$ws.Range("H88").Value = 515.625
$ws.Range("I88").Value = 680
$ws.Range("J88").Value = 241.66667
$ws.Range("K88").Value = 680
$ws.Range("L88").Value = 241.66667
$ws.Range("M88").Value = -274
$ws.Range("N88").Value = -1053.66667

$ws.Range("H91").Value = 515.625
$ws.Range("I91").Value = 680
$ws.Range("J91").Value = 241.66667
$ws.Range("K91").Value = 680
$ws.Range("L91").Value = 241.66667
$ws.Range("M91").Value = 724
$ws.Range("N91").Value = -3049.66667

$ws.Range("H94").Value = 50000
$ws.Range("J94").Value = 50000
$ws.Range("L94").Value = 50000
$ws.Range("N94").Value = -51802

$ws.Range("H122").Value = 478438.75
$ws.Range("I122").Value = 770638.0600000001
$ws.Range("J122").Value = 3614.875
$ws.Range("K122").Value = 2311914.18
$ws.Range("L122").Value = 10844.625
$ws.Range("M122").Value = -2309464.18
$ws.Range("N122").Value = -15744.625

$ws.Range("H132").Value = 4771.4287
$ws.Range("I132").Value = 3480
$ws.Range("K132").Value = 10440
$ws.Range("M132").Value = -7910

$ws.Range("H136").Value = 3305.6667
$ws.Range("I136").Value = 3305.6667
$ws.Range("K136").Value = 9917.000100000001
$ws.Range("M136").Value = -7367.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3442.9285
$ws.Range("I86").Value = 3253.875
$ws.Range("J86").Value = 3695
$ws.Range("K86").Value = 3253.875
$ws.Range("L86").Value = 3695
$ws.Range("M86").Value = -2130.875
$ws.Range("N86").Value = -5941

$ws.Range("H89").Value = 3442.9285
$ws.Range("I89").Value = 3253.875
$ws.Range("J89").Value = 3695
$ws.Range("K89").Value = 16269.375
$ws.Range("L89").Value = 18475
$ws.Range("M89").Value = -10653.375
$ws.Range("N89").Value = -29707

$ws.Range("H105").Value = 4921.1577
$ws.Range("I105").Value = 4018.0715
$ws.Range("J105").Value = 7449.8
$ws.Range("K105").Value = 4018.0715
$ws.Range("L105").Value = 7449.8
$ws.Range("M105").Value = -2271.0715
$ws.Range("N105").Value = -10943.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1033.25
$ws.Range("I16").Value = 923.7143
$ws.Range("K16").Value = 923.7143
$ws.Range("M16").Value = -636.7143

$ws.Range("H92").Value = 39867
$ws.Range("J92").Value = 39867
$ws.Range("L92").Value = 39867
$ws.Range("N92").Value = -44859

$ws.Range("H107").Value = 1600
$ws.Range("I107").Value = 1250
$ws.Range("K107").Value = 1250
$ws.Range("M107").Value = 670

$ws.Range("H113").Value = 1033.25
$ws.Range("I113").Value = 923.7143
$ws.Range("K113").Value = 923.7143
$ws.Range("M113").Value = 1246.2857

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1019.3333
$ws.Range("I14").Value = 1019.3333
$ws.Range("K14").Value = 3057.9999
$ws.Range("M14").Value = -2884.9999

$ws.Range("H131").Value = 2510
$ws.Range("I131").Value = 1364.4445
$ws.Range("J131").Value = 3982.8572
$ws.Range("K131").Value = 4093.3335
$ws.Range("L131").Value = 11948.5716
$ws.Range("M131").Value = 946.6664999999998
$ws.Range("N131").Value = -22028.5716

$ws.Range("H138").Value = 2932.375
$ws.Range("I138").Value = 1238
$ws.Range("J138").Value = 8015.5
$ws.Range("K138").Value = 3714
$ws.Range("L138").Value = 24046.5
$ws.Range("M138").Value = 1426
$ws.Range("N138").Value = -34326.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1307.7667
$ws.Range("I102").Value = 426.375
$ws.Range("K102").Value = 426.375
$ws.Range("M102").Value = 1195.625

$ws.Range("H107").Value = 795.8261
$ws.Range("I107").Value = 573.2727
$ws.Range("K107").Value = 573.2727
$ws.Range("M107").Value = 1346.7273

$ws.Range("H122").Value = 36592.766
$ws.Range("J122").Value = 94740.37
$ws.Range("L122").Value = 284221.11
$ws.Range("N122").Value = -289121.11

$ws.Range("H132").Value = 3088.7273
$ws.Range("I132").Value = 2169.7144
$ws.Range("K132").Value = 6509.1432
$ws.Range("M132").Value = -3979.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 6708332.5
$ws.Range("I2").Value = 20000000
$ws.Range("K2").Value = 20000000
$ws.Range("M2").Value = -19999888

$ws.Range("H46").Value = 2969.1738
$ws.Range("I46").Value = 2120.6428
$ws.Range("J46").Value = 4289.1113
$ws.Range("K46").Value = 2120.6428
$ws.Range("L46").Value = 4289.1113
$ws.Range("M46").Value = -1932.6428
$ws.Range("N46").Value = -4665.1113

$ws.Range("H136").Value = 3538.6
$ws.Range("I136").Value = 3587
$ws.Range("K136").Value = 10761
$ws.Range("M136").Value = -8211

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 499
$ws.Range("I2").Value = 499
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 499
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -387
$ws.Range("N2").ClearContents()

$ws.Range("H132").Value = 1553.1111
$ws.Range("I132").Value = 1636.8
$ws.Range("K132").Value = 4910.4
$ws.Range("M132").Value = -2380.4

$ws.Range("H136").Value = 3061.8125
$ws.Range("J136").Value = 8471
$ws.Range("L136").Value = 25413
$ws.Range("N136").Value = -30513

Write-Output "Applied $([int]224) cell updates across 8 sheets"
